$d = $word.ActiveDocument

$old = "same directory as the programme script."
$new = "same directory as the programme script and need to remain in the folder for the programme output to be valid. For the programme to recognise and select files, they must contain a ‘-’. This was based on the naming conventions of the data used to design the programme and was maintained as the best practice for file naming includes multiple fields separated by dashes to identify the file without having to open it. This encourages the user to adopt a logical file naming structure which will be beneficial in the long run."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Output "Replace result: $found"
